# Repull data, push all data, mean calculation
# Update the dSF column (F) values that changed after recalculating
# the mean/differences for specific rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    16 = 2
    24 = 0
    29 = 0
    32 = -2
    42 = -1
    48 = -4
    54 = 2
    56 = 0
    68 = -1
    72 = -1
    77 = -3
    79 = -6
    81 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
